# Update "想去人数" (interest counts) in the F column across sheets,
# reflecting refreshed scrape totals (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1255
$ws.Range("F5").Value = 930
$ws.Range("F6").Value = 1676
$ws.Range("F11").Value = 247
$ws.Range("F12").Value = 14
$ws.Range("F14").Value = 611
$ws.Range("F15").Value = 118
$ws.Range("F16").Value = 72
$ws.Range("F20").Value = 71
$ws.Range("F21").Value = 630
$ws.Range("F25").Value = 828
$ws.Range("F26").Value = 286
$ws.Range("F28").Value = 19
$ws.Range("F30").Value = 4

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 610

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1255
$ws.Range("F6").Value = 930
$ws.Range("F7").Value = 1676
$ws.Range("F13").Value = 247
$ws.Range("F14").Value = 14
$ws.Range("F16").Value = 611
$ws.Range("F17").Value = 118
$ws.Range("F18").Value = 72
$ws.Range("F28").Value = 71
$ws.Range("F29").Value = 630
$ws.Range("F33").Value = 828
$ws.Range("F34").Value = 286
$ws.Range("F37").Value = 19
$ws.Range("F39").Value = 610
$ws.Range("F42").Value = 4

